$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 51
$ws.Range("A51").Value = 130800369
$ws.Range("B51").Value = 79243
$ws.Range("D51").Value = "NT"
$ws.Range("E51").Value = 6425
$ws.Range("F51").Value = "Garnlav"
$ws.Range("G51").Value = "Alectoria sarmentosa"
$ws.Range("H51").Value = "(Ach.) Ach."
$ws.Range("I51").Value = "80"
$ws.Range("J51").Value = "bålar"
$ws.Range("Q51").Value = 444704
$ws.Range("R51").Value = 7025567
$ws.Range("Z51").Value = "13:25"
$ws.Range("AB51").Value = "13:25"
$ws.Range("AC51").Value = "80 bålar på gammal levande tall (ca 200 år, 23 cm dbh) i gles gammal barrblandskog, längsta bål 40 cm"
$ws.Range("AJ51").Value = "tall"
$ws.Range("AK51").Value = "Pinus sylvestris"
$ws.Range("AO51").Value = "Pinus sylvestris"

# Row 52
$ws.Range("A52").Value = 130800359
$ws.Range("B52").Value = 75221
$ws.Range("D52").Value = "LC"
$ws.Range("E52").Value = 6428
$ws.Range("F52").Value = "Rostfläck"
$ws.Range("G52").Value = "Arthonia vinosa"
$ws.Range("H52").Value = "Leight."
$ws.Range("I52").Value = ""
$ws.Range("J52").Value = ""
$ws.Range("Q52").Value = 444731
$ws.Range("R52").Value = 7025516
$ws.Range("Z52").Value = "13:54"
$ws.Range("AB52").Value = "13:54"
$ws.Range("AC52").Value = "Vid basen av grov gammal levande sälg (50 cm dbh) i gammal granskog"
$ws.Range("AJ52").Value = "sälg"
$ws.Range("AK52").Value = "Salix caprea"
$ws.Range("AO52").Value = "Salix caprea"

# Row 54
$ws.Range("A54").Value = 130800375
$ws.Range("B54").Value = 83214
$ws.Range("D54").Value = "VU"
$ws.Range("E54").Value = 492
$ws.Range("F54").Value = "Smalskaftslav"
$ws.Range("G54").Value = "Chaenotheca gracilenta"
$ws.Range("H54").Value = "(Ach.) J.Mattsson & Middelb."
$ws.Range("J54").Value = ""
$ws.Range("K54").Value = ""
$ws.Range("N54").Value = ""
$ws.Range("Q54").Value = 444662
$ws.Range("R54").Value = 7025556
$ws.Range("Z54").Value = "13:09"
$ws.Range("AB54").Value = "13:09"
$ws.Range("AC54").Value = "Vid basen av grov björkhögstubbe (30 cm dbh) i gammal granskog med inslag av tallöverståndare"
$ws.Range("AF54").Value = ""
$ws.Range("AJ54").Value = "glasbjörk"
$ws.Range("AK54").Value = "Betula pubescens"
$ws.Range("AO54").Value = "Betula pubescens"

# Row 55
$ws.Range("A55").Value = 130800352
$ws.Range("B55").Value = 79243
$ws.Range("D55").Value = "NT"
$ws.Range("E55").Value = 6425
$ws.Range("F55").Value = "Garnlav"
$ws.Range("G55").Value = "Alectoria sarmentosa"
$ws.Range("H55").Value = "(Ach.) Ach."
$ws.Range("J55").Value = ""
$ws.Range("K55").Value = ""
$ws.Range("N55").Value = ""
$ws.Range("Q55").Value = 444700
$ws.Range("R55").Value = 7025517
$ws.Range("Z55").Value = "14:43"
$ws.Range("AB55").Value = "14:43"
$ws.Range("AC55").Value = "Rikligt på gammal levande gran i gammal granskog"
$ws.Range("AF55").Value = ""
$ws.Range("AJ55").Value = "gran"
$ws.Range("AK55").Value = "Picea abies"
$ws.Range("AO55").Value = "Picea abies"

# Row 56
$ws.Range("A56").Value = 130800364
$ws.Range("B56").Value = 80384
$ws.Range("D56").Value = "LC"
$ws.Range("E56").Value = 6464
$ws.Range("F56").Value = "Luddlav"
$ws.Range("G56").Value = "Nephroma resupinatum"
$ws.Range("H56").Value = "(L.) Ach."
$ws.Range("Q56").Value = 444717
$ws.Range("R56").Value = 7025526
$ws.Range("Z56").Value = "13:48"
$ws.Range("AB56").Value = "13:48"
$ws.Range("AC56").Value = "På bark på stam av levande lutande gammal sälg i gammal granskog"

# Row 57
$ws.Range("A57").Value = 130800354
$ws.Range("B57").Value = 79243
$ws.Range("D57").Value = "NT"
$ws.Range("E57").Value = 6425
$ws.Range("F57").Value = "Garnlav"
$ws.Range("G57").Value = "Alectoria sarmentosa"
$ws.Range("H57").Value = "(Ach.) Ach."
$ws.Range("Q57").Value = 444716
$ws.Range("R57").Value = 7025439
$ws.Range("Z57").Value = "14:16"
$ws.Range("AB57").Value = "14:16"
$ws.Range("AC57").Value = "På gammal gran i gammal granskog"
$ws.Range("AJ57").Value = "gran"
$ws.Range("AK57").Value = "Picea abies"
$ws.Range("AO57").Value = "Picea abies"

# Row 58
$ws.Range("A58").Value = 130800361
$ws.Range("B58").Value = 75333
$ws.Range("E58").Value = 1460
$ws.Range("F58").Value = "Rosa skärelav"
$ws.Range("G58").Value = "Schismatomma pericleum"
$ws.Range("Q58").Value = 444731
$ws.Range("R58").Value = 7025516
$ws.Range("AB58").Value = "13:54"
$ws.Range("AC58").Value = "Vid basen av grov gammal levande sälg (50 cm dbh) i gammal granskog"
$ws.Range("AJ58").Value = "sälg"
$ws.Range("AK58").Value = "Salix caprea"
$ws.Range("AO58").Value = "Salix caprea"
